$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43:71 down to 44:72
$ws.Rows(43).Insert()

# Populate the new row 43 with the new weekly record
$ws.Cells.Item(43, 1).Value = 4
$ws.Cells.Item(43, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(43, 3).Value = "Los Lagos"
$ws.Cells.Item(43, 4).Value = 45216
$ws.Cells.Item(43, 5).Value = 10
$ws.Cells.Item(43, 6).Value = "Fruta"
$ws.Cells.Item(43, 7).Value = 100101
$ws.Cells.Item(43, 8).Value = "Berries"
$ws.Cells.Item(43, 9).Value = 100101001
$ws.Cells.Item(43, 10).Value = "Arándano (blue)"
$ws.Cells.Item(43, 11).Value = "Sin especificar"
$ws.Cells.Item(43, 12).Value = "Primera"
$ws.Cells.Item(43, 13).Value = 100
$ws.Cells.Item(43, 14).Value = 12000
$ws.Cells.Item(43, 15).Value = 12000
$ws.Cells.Item(43, 16).Value = 12000
$ws.Cells.Item(43, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(43, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(43, 19).Value = 6000
$ws.Cells.Item(43, 20).Value = 2
